$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Remove the fill style that was applied to B17 (revert it back to default/no explicit style)
$ws.Range("B17").ClearFormats()

# Update C17's value from 123456 to 12345 (numeric)
$ws.Range("C17").Value = 12345

# Update the active selection on the sheet to E15
$ws.Activate()
$ws.Range("E15").Select()
